$wb = $excel.ActiveWorkbook

$wsRegistrants = $wb.Worksheets.Item("registrants")
$wsServices = $wb.Worksheets.Item("services")

# --- sheet "registrants" (sheet1) ---
# D6 held the shared string "PT_BR" which changes (in the shared-strings table)
# to "IT_IT" everywhere it is used - no other cell in this sheet references
# that string, so updating D6's text covers it.
$wsRegistrants.Range("D6").Value = "IT_IT"

# Update the selection on sheet1 to match the new view.
$wsRegistrants.Range("B11").Select()

# --- sheet "services" (sheet2) ---
$wsServices.Range("D3").Value = "DE_DE"
$wsServices.Range("D5").Value = "EN_GB"
$wsServices.Range("D6").Value = "EN_GB"
$wsServices.Range("D7").Value = "FR_FR"

# Update the selection on sheet2 to match the new view.
$wsServices.Range("D2").Select()
